$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows to grow the table from 4 data rows to 6 data rows
# (old row2 -> row3, old row3 -> row4, old row4 -> row5, old row5 -> row6 after first insert;
#  then a second insert at row4 pushes row4.. down once more)
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(4).Insert()

# Make sure the whole data range is formatted as text so numeric-looking strings stay text
$ws.Range("A2:T7").NumberFormat = "@"

# Row 2: /model.json
$ws.Range("A2").Value = "/model.json"
$ws.Range("B2").Value = "12.0"
$ws.Range("C2").Value = "0.0"
$ws.Range("D2").Value = "6.0"
$ws.Range("E2").Value = "-6.0"
$ws.Range("F2").Value = "2.25"
$ws.Range("G2").Value = "-0.375"
$ws.Range("H2").Value = "8.0"
$ws.Range("I2").Value = "5.84375"
$ws.Range("J2").Value = "0.0"
$ws.Range("K2").Value = "0.1953125"
$ws.Range("L2").Value = "59.77620442708333"
$ws.Range("M2").Value = "-0.0"
$ws.Range("N2").Value = "2.490234375"
$ws.Range("O2").Value = "-0.0"
$ws.Range("P2").Value = "101.23909505208337"
$ws.Range("Q2").Value = "-0.0"
$ws.Range("R2").Value = "2.490234375"
$ws.Range("S2").Value = "-0.0"
$ws.Range("T2").Value = "158.47981770833337"

# Row 3: /stab/stab.json
$ws.Range("A3").Value = "/stab/stab.json"
$ws.Range("B3").Value = "2.0"
$ws.Range("C3").Value = "0.0"
$ws.Range("D3").Value = "3.0"
$ws.Range("E3").Value = "-3.0"
$ws.Range("F3").Value = "0.125"
$ws.Range("G3").Value = "0.0"
$ws.Range("H3").Value = "1.5"
$ws.Range("I3").Value = "1.0"
$ws.Range("J3").Value = "0.0"
$ws.Range("K3").Value = "0.0625"
$ws.Range("L3").Value = "4.501953125"
$ws.Range("M3").Value = "-0.0"
$ws.Range("N3").Value = "-0.0"
$ws.Range("O3").Value = "-0.0"
$ws.Range("P3").Value = "0.501953125"
$ws.Range("Q3").Value = "-0.0"
$ws.Range("R3").Value = "-0.0"
$ws.Range("S3").Value = "-0.0"
$ws.Range("T3").Value = "5.0"

# Row 4: /body/body.json
$ws.Range("A4").Value = "/body/body.json"
$ws.Range("B4").Value = "12.0"
$ws.Range("C4").Value = "0.0"
$ws.Range("D4").Value = "0.125"
$ws.Range("E4").Value = "-0.125"
$ws.Range("F4").Value = "2.25"
$ws.Range("G4").Value = "-0.25"
$ws.Range("H4").Value = "2.0"
$ws.Range("I4").Value = "7.25"
$ws.Range("J4").Value = "0.0"
$ws.Range("K4").Value = "0.3125"
$ws.Range("L4").Value = "0.7923177083333333"
$ws.Range("M4").Value = "-0.0"
$ws.Range("N4").Value = "-2.34375"
$ws.Range("O4").Value = "-0.0"
$ws.Range("P4").Value = "28.325520833333314"
$ws.Range("Q4").Value = "-0.0"
$ws.Range("R4").Value = "-2.34375"
$ws.Range("S4").Value = "-0.0"
$ws.Range("T4").Value = "27.550130208333314"

# Row 5: /body/rudder/rudder.json
$ws.Range("A5").Value = "/body/rudder/rudder.json"
$ws.Range("B5").Value = "2.0"
$ws.Range("C5").Value = "0.0"
$ws.Range("D5").Value = "0.0625"
$ws.Range("E5").Value = "-0.0625"
$ws.Range("F5").Value = "2.0"
$ws.Range("G5").Value = "0.0"
$ws.Range("H5").Value = "0.5"
$ws.Range("I5").Value = "1.0"
$ws.Range("J5").Value = "0.0"
$ws.Range("K5").Value = "1.0"
$ws.Range("L5").Value = "0.16731770833333326"
$ws.Range("M5").Value = "-0.0"
$ws.Range("N5").Value = "-0.0"
$ws.Range("O5").Value = "-0.0"
$ws.Range("P5").Value = "0.33333333333333326"
$ws.Range("Q5").Value = "-0.0"
$ws.Range("R5").Value = "-0.0"
$ws.Range("S5").Value = "-0.0"
$ws.Range("T5").Value = "0.16731770833333326"

# Row 6: /body/stick/stick.json
$ws.Range("A6").Value = "/body/stick/stick.json"
$ws.Range("B6").Value = "12.0"
$ws.Range("C6").Value = "0.0"
$ws.Range("D6").Value = "0.125"
$ws.Range("E6").Value = "-0.125"
$ws.Range("F6").Value = "0.25"
$ws.Range("G6").Value = "-0.25"
$ws.Range("H6").Value = "1.5"
$ws.Range("I6").Value = "6.0"
$ws.Range("J6").Value = "0.0"
$ws.Range("K6").Value = "0.0"
$ws.Range("L6").Value = "0.0390625"
$ws.Range("M6").Value = "-0.0"
$ws.Range("N6").Value = "-0.0"
$ws.Range("O6").Value = "-0.0"
$ws.Range("P6").Value = "18.03125"
$ws.Range("Q6").Value = "-0.0"
$ws.Range("R6").Value = "-0.0"
$ws.Range("S6").Value = "-0.0"
$ws.Range("T6").Value = "18.0078125"

# Row 7: /wing/wing.json
$ws.Range("A7").Value = "/wing/wing.json"
$ws.Range("B7").Value = "3.0"
$ws.Range("C7").Value = "0.0"
$ws.Range("D7").Value = "6.0"
$ws.Range("E7").Value = "-6.0"
$ws.Range("F7").Value = "0.125"
$ws.Range("G7").Value = "0.0"
$ws.Range("H7").Value = "4.5"
$ws.Range("I7").Value = "1.5"
$ws.Range("J7").Value = "0.0"
$ws.Range("K7").Value = "0.0625"
$ws.Range("L7").Value = "54.005859375"
$ws.Range("M7").Value = "-0.0"
$ws.Range("N7").Value = "-0.0"
$ws.Range("O7").Value = "-0.0"
$ws.Range("P7").Value = "3.380859375"
$ws.Range("Q7").Value = "-0.0"
$ws.Range("R7").Value = "-0.0"
$ws.Range("S7").Value = "-0.0"
$ws.Range("T7").Value = "57.375"
